# testDataUS.xlsx - refresh the sample New Hire test-data row used by the
# automation tests (commit: "Mad the changes in webActionPage& created
# sepreate pages for JobDetalsand Propose Compensation,etc under
# commonpage folder.") - the sample row's person/grade/ID fields were
# updated to a new test fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# GivenName / FamilyName (columns E/F): Dennis Ruecker -> Tara Robel
$ws.Range("E2").Value = "Tara"
$ws.Range("F2").Value = "Robel"

# Step (column AM): "Initial Step - CGM" -> "N/A"
$ws.Range("AM2").Value = "N/A"

# AddEditID1 / national ID value (column AP): new generated ID
$ws.Range("AP2").Value = 968221814
